# Lesson 3 ("Vistas y Validaciones") spacing paragraph rework.
#
# 1) Insert two new empty spacer paragraphs right after the "Vistas" link
#    paragraph (before the existing empty ind=1440 paragraph that precedes
#    "Taller 1"): one plain spacer (no style/indent) and one indented
#    (left=720) spacer.
# 2) Collapse the pair of empty paragraphs that sits between "Taller 1" and
#    "Validaciones de Documentos" into a single plain spacer paragraph
#    (drop the ListParagraph style + 1440 indent).
# 3) Remove the extra empty ListParagraph-styled paragraph that follows
#    "Taller 2".

$d = $word.ActiveDocument

function Find-ParagraphIndexByText($searchText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $searchText"
    }
    $target = $rng.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($target -ge $cand.Range.Start -and $target -lt $cand.Range.End) {
            return $i
        }
    }
    throw "Paragraph containing text not found: $searchText"
}

# Package-wrapped WordOpenXML used to overwrite a single paragraph's
# contents via Range.InsertXML (replaces whatever the range currently
# holds with exactly this <w:p>).
function Make-ParaPkg($innerPPr) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p><w:pPr>' + $innerPPr + '</w:pPr></w:p></w:body>' + `
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$rPrArial = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$spacingAuto = '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>'

$plainSpacerPkg = Make-ParaPkg ($spacingAuto + $rPrArial)
$indentedSpacerPkg = Make-ParaPkg ($spacingAuto + '<w:ind w:left="720"/>' + $rPrArial)

# "Validaciones de Documentos" is unique in the document and anchors the
# whole Lesson-3 block; every other target paragraph sits at a fixed
# offset relative to it.
$A = Find-ParagraphIndexByText("Validaciones de Documentos")

# Mutate starting from the highest paragraph index downward so earlier
# (lower-numbered) paragraphs we still need to address keep their index
# valid.

# --- (3) delete the empty ListParagraph paragraph after "Taller 2" ---
$d.Paragraphs.Item($A + 3).Range.Delete()

# --- (2) merge the two empty paragraphs after "Taller 1" into one ---
# Strip the pPr of the second (ind=1440) paragraph down to just spacing+rPr,
# then delete the first (plain ListParagraph) paragraph entirely.
$d.Paragraphs.Item($A - 1).Range.InsertXML($plainSpacerPkg)
$d.Paragraphs.Item($A - 2).Range.Delete()

# --- (1) insert two new spacer paragraphs before the paragraph preceding
#          "Taller 1" (the empty ind=1440 paragraph right after "Vistas") ---
$precedingTaller1 = $d.Paragraphs.Item($A - 4)
$precedingTaller1.Range.InsertParagraphBefore()
$d.Paragraphs.Item($A - 4).Range.InsertXML($plainSpacerPkg)

$precedingTaller1Again = $d.Paragraphs.Item($A - 3)
$precedingTaller1Again.Range.InsertParagraphBefore()
$d.Paragraphs.Item($A - 3).Range.InsertXML($indentedSpacerPkg)
